# DISK section ("free", "total", "inodes_free" rows): the per-metric Notes
# column previously said the disk metrics required host-level monitoring.
# Replace that note with the same "cgroups can't monitor this directly"
# wording already used elsewhere in the sheet (e.g. usage_iowait / io_time /
# io_util), which also makes the old "Requires host-level monitoring" string
# unused/removed from the workbook's shared string table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$note = "cgroups não consegue monitorar isso diretamente"
$ws.Range("C15").Value = $note
$ws.Range("C16").Value = $note
$ws.Range("C17").Value = $note

# Leave the sheet's active cell/selection where the author last left it.
$ws.Range("C24").Select()
